$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply all cell updates from the diff. D and E columns use a text-format
# + style-reset pattern so that numeric-looking strings (e.g. "210.09")
# are stored as text (inline/shared string), matching the source data,
# rather than being auto-coerced into numeric cell values by Excel.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "26.335.31"
Set-TextValue "E2" "  -1.09%  "
Set-TextValue "D3" "1.588.43"
Set-TextValue "E3" "  -0.61%  "
Set-TextValue "D5" "210.09"
Set-TextValue "E6" "  -1.24%  "
Set-TextValue "E8" "  -0.76%  "
Set-TextValue "E9" "  -0.50%  "
Set-TextValue "D10" "19.49"
Set-TextValue "E10" "  -0.37%  "
Set-TextValue "D11" "0.0846"
Set-TextValue "E11" "  +0.04%  "
Set-TextValue "D12" "1.813.13"
Set-TextValue "E12" "  -0.53%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.590.98"
Set-TextValue "E13" "  -0.44%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D14" "4.07"
Set-TextValue "E14" "  +0.69%  "
Set-TextValue "E15" "  -0.94%  "
Set-TextValue "D16" "64.31"
Set-TextValue "E16" "  -0.44%  "
Set-TextValue "D17" "26.343.02"
Set-TextValue "E17" "  -1.01%  "
Set-TextValue "E18" "  -1.58%  "
Set-TextValue "D19" "7.46"
Set-TextValue "E19" "  +5.16%  "
Set-TextValue "D20" "210.74"
Set-TextValue "E20" "  +1.14%  "
Set-TextValue "E21" "  -0.39%  "
Set-TextValue "E22" "  -0.19%  "
Set-TextValue "D23" "2.15"
Set-TextValue "E23" "  -3.89%  "
Set-TextValue "D24" "8.92"
Set-TextValue "E24" "  -0.24%  "
Set-TextValue "D25" "145.03"
Set-TextValue "E25" "  +0.87%  "
Set-TextValue "E26" "  -0.33%  "
Set-TextValue "D27" "7.06"
Set-TextValue "E27" "  -1.22%  "
Set-TextValue "E28" "  -0.48%  "
Set-TextValue "D29" "15.23"
Set-TextValue "E29" "  -0.28%  "
Set-TextValue "E30" "  -0.45%  "
Set-TextValue "E31" "  -0.29%  "
Set-TextValue "E32" "  -0.81%  "
Set-TextValue "E33" "  +1.12%  "
Set-TextValue "D34" "1.304.03"
Set-TextValue "E34" "  +1.99%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D35" "2.44"
Set-TextValue "E35" "  -1.71%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D36" "0.611"
Set-TextValue "E36" "  +2.20%  "
Set-TextValue "E37" "  -0.87%  "
Set-TextValue "E38" "  +0.31%  "
Set-TextValue "E39" "  -13.27%  "
Set-TextValue "D40" "0.809"
Set-TextValue "E40" "  -1.69%  "
Set-TextValue "E41" "  -0.36%  "
Set-TextValue "E42" "  +3.61%  "
Set-TextValue "D43" "0.768"
Set-TextValue "E43" "  -0.60%  "
Set-TextValue "E44" "  -1.50%  "
Set-TextValue "D45" "62.48"
Set-TextValue "E45" "  -0.18%  "
Set-TextValue "D46" "1.724.87"
Set-TextValue "E46" "  -0.52%  "
Set-TextValue "D47" "87.74"
Set-TextValue "E47" "  -2.09%  "
Set-TextValue "E48" "  -5.05%  "
Set-TextValue "E49" "  -1.45%  "
Set-TextValue "D50" "0.0979"
Set-TextValue "E50" "  -4.45%  "
